$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column in H, matching the style of the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Populate the save-flag values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
